# Weekly update: prepend a new day's price record for Espárragos at
# Macroferia Regional de Talca. A new row is inserted at row 4 (pushing the
# existing rows 4-84 down to 5-85) and populated with the latest reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 4, shifting old rows 4..84 -> 5..85
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new record
$ws.Range("A4").Value() = 5
$ws.Range("B4").Value() = "Macroferia Regional de Talca"
$ws.Range("C4").Value() = "Maule"
$ws.Range("D4").Value() = 44882
$ws.Range("E4").Value() = 7
$ws.Range("F4").Value() = 300000000
$ws.Range("G4").Value() = "Espárragos"
$ws.Range("H4").Value() = "Sin especificar"
$ws.Range("I4").Value() = "Primera"
$ws.Range("J4").Value() = 3000
$ws.Range("K4").Value() = 1000
$ws.Range("L4").Value() = 1000
$ws.Range("M4").Value() = 1000
$ws.Range("N4").Value() = "$/kilo"
$ws.Range("O4").Value() = "Región del Maule"
$ws.Range("P4").Value() = 1000
$ws.Range("Q4").Value() = 1
$ws.Range("R4").Value() = "Hortaliza"
